$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 4): E4 "Dependency", F4 "Description" ---
# These need the same bold header style already used by A4:D4 (fontId1:
# bold, 14pt Calibri). Copy A4's formatting only (not its value) onto the
# new header cells so the existing style is reused instead of a new one
# being allocated.
$ws.Range("E4").Value = "Dependency"
$ws.Range("F4").Value = "Description"
$ws.Range("A4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null

# --- Existing valve rows 5-8 gain a new "#" column (E) ---
$ws.Range("E5").Value = "#"
$ws.Range("E6").Value = "#"
$ws.Range("E7").Value = "#"
$ws.Range("E8").Value = "#"

# --- New valve rows 9-10 ---
$ws.Range("A9").Value = "Valve"
$ws.Range("B9").Value = "Tap-Off-A Valve"
$ws.Range("C9").Value = "Feed System"
$ws.Range("D9").Value = "N"
$ws.Range("E9").Value = "#"

$ws.Range("A10").Value = "Valve "
$ws.Range("B10").Value = "Tap-Off-B Valve"
$ws.Range("C10").Value = "Feed System"
$ws.Range("D10").Value = "N"
$ws.Range("E10").Value = "#"

# --- New pressure-transducer rows 12-15 (row 11 left blank) ---
$ws.Range("A12").Value = "Pressure Transducer"
$ws.Range("B12").Value = "Post-Inlet-Fuel transducer"
$ws.Range("C12").Value = "Avionics"
$ws.Range("D12").Value = "R"
$ws.Range("E12").Value = "Fuel Flow in propellant line / Mass Flow"

$ws.Range("A13").Value = "Pressure Transducer"
$ws.Range("B13").Value = "Post-Inlet-LOX transducer"
$ws.Range("C13").Value = "Avionics"
$ws.Range("D13").Value = "R"
$ws.Range("E13").Value = "LOX Flow in propellant line /  Mass Flow"

$ws.Range("A14").Value = "Pressure Transducer"
$ws.Range("B14").Value = "Pre-Chamber Fuel transducer"
$ws.Range("C14").Value = "Avionics"
$ws.Range("D14").Value = "R"
$ws.Range("E14").Value = "Fuel Flow in propellant line /  Mass Flow"

$ws.Range("A15").Value = "Pressure Transducer"
$ws.Range("B15").Value = "Pre-Chamber LOX transducer"
$ws.Range("C15").Value = "Avionics"
$ws.Range("D15").Value = "R"
$ws.Range("E15").Value = "LOX Flow in propellant line /  Mass Flow"

# --- New thermocouple rows 17-18 (row 16 left blank) ---
$ws.Range("A17").Value = "Thermocouple"
$ws.Range("B17").Value = "Pre-Chamber Fuel thermocouple"
$ws.Range("C17").Value = "Avionics"
$ws.Range("D17").Value = "R"
$ws.Range("E17").Value = "Temp in Fuel line"

$ws.Range("A18").Value = "Thermocouple"
$ws.Range("B18").Value = "Pre-Chamber LOX thermocouple"
$ws.Range("C18").Value = "Avionics"
$ws.Range("D18").Value = "R"
$ws.Range("E18").Value = "Temp in LOX line"

# --- Column widths (A-F) ---
$ws.Columns.Item(1).ColumnWidth = 24.666666666666668
$ws.Columns.Item(2).ColumnWidth = 27.998697916666668
$ws.Columns.Item(3).ColumnWidth = 20.166666666666668
$ws.Columns.Item(4).ColumnWidth = 30.166666666666668
$ws.Columns.Item(5).ColumnWidth = 33.998697916666664
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666

# --- Selection moves to E17, matching the saved cursor state ---
$ws.Range("E17").Select() | Out-Null
